$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '46.113.66'
$ws.Cells.Item(2, 5).Value = '  -1.96%  '
$ws.Cells.Item(3, 4).Value = '2.334.56'
$ws.Cells.Item(3, 5).Value = '  -0.03%  '
$ws.Cells.Item(4, 5).Value = '  +0.20%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '300.85'
$ws.Cells.Item(5, 5).Value = '  -1.70%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '98.40'
$ws.Cells.Item(6, 5).Value = '  +0.04%  '
$ws.Cells.Item(7, 5).Value = '  -1.71%  '
$ws.Cells.Item(8, 5).Value = '  +0.11%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.510'
$ws.Cells.Item(9, 5).Value = '  -5.64%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '34.63'
$ws.Cells.Item(10, 5).Value = '  -4.14%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0790'
$ws.Cells.Item(11, 5).Value = '  -3.28%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '7.11'
$ws.Cells.Item(12, 5).Value = '  -5.20%  '
$ws.Cells.Item(13, 5).Value = '  -1.52%  '
$ws.Cells.Item(14, 4).Value = '2.694.81'
$ws.Cells.Item(14, 5).Value = '  +0.17%  '
$ws.Cells.Item(15, 4).Value = '2.334.05'
$ws.Cells.Item(15, 5).Value = '  -0.14%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '13.65'
$ws.Cells.Item(16, 5).Value = '  -3.80%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.804'
$ws.Cells.Item(17, 5).Value = '  -4.40%  '
$ws.Cells.Item(18, 4).Value = '46.086.74'
$ws.Cells.Item(18, 5).Value = '  -1.68%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '12.66'
$ws.Cells.Item(19, 5).Value = '  -7.31%  '
$ws.Cells.Item(20, 4).Value = '0.0₃0966'
$ws.Cells.Item(20, 5).Value = '  +1.05%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '5.96'
$ws.Cells.Item(21, 5).Value = '  -4.18%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '66.50'
$ws.Cells.Item(22, 5).Value = '  -2.43%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '243.99'
$ws.Cells.Item(23, 5).Value = '  -4.01%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '2.81'
$ws.Cells.Item(24, 5).Value = '  -5.79%  '
$ws.Cells.Item(25, 5).Value = '  +0.00%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '1.90'
$ws.Cells.Item(26, 5).Value = '  -5.58%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '40.26'
$ws.Cells.Item(27, 5).Value = '  -5.15%  '
$ws.Cells.Item(28, 5).Value = '  -2.76%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '9.67'
$ws.Cells.Item(29, 5).Value = '  -2.88%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '20.73'
$ws.Cells.Item(30, 5).Value = '  +1.71%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '3.60'
$ws.Cells.Item(31, 5).Value = '  +14.59%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '2.81'
$ws.Cells.Item(32, 5).Value = '  +6.64%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '5.43'
$ws.Cells.Item(33, 5).Value = '  -7.07%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '144.65'
$ws.Cells.Item(34, 5).Value = '  -1.26%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.0771'
$ws.Cells.Item(35, 5).Value = '  -6.00%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.112'
$ws.Cells.Item(36, 5).Value = '  -3.15%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.116'
$ws.Cells.Item(37, 5).Value = '  -3.07%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '1.79'
$ws.Cells.Item(38, 5).Value = '  -2.10%  '
$ws.Cells.Item(39, 5).Value = '  +7.45%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '3.87'
$ws.Cells.Item(40, 5).Value = '  -3.69%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.0298'
$ws.Cells.Item(41, 5).Value = '  -4.65%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '3.18'
$ws.Cells.Item(42, 5).Value = '  -6.59%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.999'
$ws.Cells.Item(43, 5).Value = '  +0.08%  '
$ws.Cells.Item(44, 4).Value = '1.858.00'
$ws.Cells.Item(44, 5).Value = '  +2.99%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '90.27'
$ws.Cells.Item(45, 5).Value = '  -2.22%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '1.81'
$ws.Cells.Item(46, 5).Value = '  -8.91%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.185'
$ws.Cells.Item(47, 5).Value = '  -5.34%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '70.03'
$ws.Cells.Item(48, 5).Value = '  -6.65%  '
$ws.Cells.Item(49, 4).Value = '2.563.21'
$ws.Cells.Item(49, 5).Value = '  -0.10%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '95.69'
$ws.Cells.Item(50, 5).Value = '  -3.55%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '4.75'
$ws.Cells.Item(51, 5).Value = '  -1.69%  '
